# Update "Cultivo" (crop) values for the rows where the author
# reclassified the crop to the newly-introduced "Trigo" (wheat) /
# "Arroz" (rice) entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 4).Value = "Trigo"
$ws.Cells.Item(10, 4).Value = "Arroz"
$ws.Cells.Item(13, 4).Value = "trigo"
$ws.Cells.Item(18, 4).Value = "Arroz"
$ws.Cells.Item(24, 4).Value = "Trigo"
$ws.Cells.Item(26, 4).Value = "Arroz"
$ws.Cells.Item(35, 4).Value = "Trigo"
$ws.Cells.Item(37, 4).Value = "Arroz"
$ws.Cells.Item(46, 4).Value = "Arroz"
$ws.Cells.Item(49, 4).Value = "Trigo"
$ws.Cells.Item(55, 4).Value = "Arroz"
$ws.Cells.Item(66, 4).Value = "Trigo"
$ws.Cells.Item(85, 4).Value = "Arroz"
$ws.Cells.Item(89, 4).Value = "Trigo"
$ws.Cells.Item(106, 4).Value = "Arroz"
$ws.Cells.Item(112, 4).Value = "Trigo"

# Correct "Mes" (month) values for a handful of rows to Marzo.
$ws.Cells.Item(85, 2).Value = "Marzo"
$ws.Cells.Item(106, 2).Value = "Marzo"
$ws.Cells.Item(112, 2).Value = "Marzo"
$ws.Cells.Item(119, 2).Value = "Marzo"

# Leave the selection where the author last left it before saving.
$ws.Range("B119").Select()
